# Apply updated cryptocurrency price/volume data to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.546.22"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.570.14"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.13%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "212.75"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.13%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "45.68"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +4.08%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "24.05"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.794.61"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "1.566.95"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "28.543.40"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("E18").Value = "  -1.33%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "229.45"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.05%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.35"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("E22").Value = "  -0.07%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "3.86"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -6.04%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.11"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("E25").Value = "  +8.81%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "151.72"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "15.02"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -1.00%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "1.392.09"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +0.68%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.54"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.36"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.56%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.61"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  -0.98%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("E42").Value = "  -0.12%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.789"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0468"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "5.50"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.00%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.970"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "62.80"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "1.706.99"
$ws.Range("E49").Value = "  -1.37%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "86.26"
$cell.Style = "Normal"
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  +0.80%  "
